$wb = $excel.ActiveWorkbook

# "Generate Report for Archive" -- the localization status moves from
# "Ready for handoff" to "In Translation" for the three language-status
# cells tracked on each sheet (Overview!E/F per-language columns, and the
# "Status" column on each per-language detail sheet). Updating the text
# makes Excel re-flow (shrink) the now-narrower "Status" columns, which is
# why their stored column widths change too.

$wsOverview = $wb.Worksheets.Item("Overview")
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsDe = $wb.Worksheets.Item("de-de")

$wsOverview.Range("E2:E4").Value = "In Translation"
$wsOverview.Range("F2:F4").Value = "In Translation"
$wsZh.Range("C2:C4").Value = "In Translation"
$wsDe.Range("C2:C4").Value = "In Translation"

# Match the narrower auto-fit width the shorter status text now produces.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZh.Columns.Item(3).ColumnWidth = 12.5
$wsDe.Columns.Item(3).ColumnWidth = 12.5
